$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("faqs-share")

# Add the new FAQ row (row 8): question, answer, include flag.
$ws.Range("A8").Value = "Can I reuse contact details for a new study?"
$ws.Range("B8").Value = "This depends on how data subjects were informed about potential reuse of their contact details: can they expect to be contacted again and for this purpose? Note that you should have obtained access to the contact details legitimately too: are you supposed to have access to their contact details in the first place? If you are uncertain about this, ask your [privacy officer](#support) for help."
$ws.Range("C8").Value = 1

# Match the recorded view state: faqs-share becomes the active sheet/tab
# with B8 selected (previously faqs-support was active).
$ws.Activate()
$ws.Range("B8").Select()
